$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.62%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.60%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.068"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.54%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05690"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.36%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.484"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.13%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8196"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.74%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8397"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.74%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.10%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06907"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.90%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02851"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.02%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09397"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.02%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001516"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.30%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04085"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-12.00%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005978"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.40%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006107"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.30%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3,761.03%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.509"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.26%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.20%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.315"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "12.64%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.94%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03195"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.37%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.77%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.563"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.91%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.68%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.38%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.003964"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-13.67%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.00009794"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03697"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.76%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005496"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "62.49%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1055"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-22.34%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001799"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-32.40%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009394"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.41%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005207"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.70%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.10%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-15.51%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002596"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.16%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"

Write-Host "Applied cryptocurrency price/volume updates"